# Apply edits to the Graph Neural Network / Inductive predictions sheet.
# 1) A batch of isolated prediction ("pred") corrections in column B.
# 2) Removal of an anomalous row (target rating "D", row 195) which shifts
#    all subsequent rows up by one and drops the final row, shrinking the
#    used range from A1:C403 to A1:C402.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Isolated column B ("pred") corrections ---
$ws.Range("B11").Value = "AAA"
$ws.Range("B15").Value = "BB"
$ws.Range("B29").Value = "BB"
$ws.Range("B30").Value = "AA"
$ws.Range("B33").Value = "BB"
$ws.Range("B39").Value = "BBB"
$ws.Range("B43").Value = "BBB"
$ws.Range("B45").Value = "BBB"
$ws.Range("B52").Value = "BB"
$ws.Range("B58").Value = "BBB"
$ws.Range("B59").Value = "B"
$ws.Range("B62").Value = "AA"
$ws.Range("B76").Value = "BB"
$ws.Range("B84").Value = "AAA"
$ws.Range("B87").Value = "A"
$ws.Range("B92").Value = "BBB"
$ws.Range("B95").Value = "BBB"
$ws.Range("B110").Value = "A"
$ws.Range("B113").Value = "AA"
$ws.Range("B116").Value = "BBB"
$ws.Range("B120").Value = "BBB"
$ws.Range("B134").Value = "BBB"
$ws.Range("B138").Value = "BBB"
$ws.Range("B140").Value = "BB"
$ws.Range("B143").Value = "A"
$ws.Range("B156").Value = "BB"
$ws.Range("B157").Value = "BB"
$ws.Range("B158").Value = "BBB"
$ws.Range("B160").Value = "AA"
$ws.Range("B161").Value = "AA"
$ws.Range("B163").Value = "A"
$ws.Range("B165").Value = "BB"
$ws.Range("B167").Value = "BBB"
$ws.Range("B169").Value = "B"
$ws.Range("B172").Value = "BB"
$ws.Range("B173").Value = "BB"
$ws.Range("B180").Value = "AA"
$ws.Range("B185").Value = "BBB"
$ws.Range("B187").Value = "A"
$ws.Range("B189").Value = "BBB"

# --- Remove the row for the "D" rated target (row 195) ---
# This shifts rows 196:403 up to 195:402 and shrinks the sheet's
# used range to A1:C402.
$ws.Rows("195:195").Delete()
